$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers
$ws.Range("A1").Value = 'Maladies'
$ws.Range("B1").Value = 'Résumés'
$ws.Range("C1").Value = 'URL'
$ws.Range("D1").Value = 'Date'
$ws.Range("E1").Value = 'Pays'

# New D1/E1 header cells need the same bold/centered header style as A1:C1
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = 'tuberculose'
$ws.Range("B2").Value = 'Tuberculose bovine : mesures renforcées face à la hausse des foyers Formation à la biosécurité, revalorisation des indemnisations, expérimentation d’un vaccin sur le blaireau, nouveau test de dépistage… Le ministère de l’Agriculture a mis en place différentes mesures de prévention pour endiguer la hausse du nombre de foyers de tuberculose bovine « dans plusieurs zones géographiques (Sud-Ouest, Normandie, Corse) et sa persistance dans d’autres (Côte d’Or, Camargue) », indique un communiqué du 5 avril. Cette maladie réglementée, transmissible à l’homme, peut être véhiculée par la faune sauvage. Les premiers tests du protocole de vaccination seront effectués en Nouvelle-Aquitaine, puis pourront être étendus « à l’ensemble des zones contaminées ». Le département, régulièrement concerné par des infections, est placé sous surveillance renforcée depuis 2015.'
$ws.Range("C2").Value = 'https://urlz.fr/lEgR'
$ws.Range("E2").Value = 'Calvados, France'

# Row 3
$ws.Range("A3").Value = 'الجلد العقدي'
$ws.Range("B3").Value = 'تعليمات عاجلة لوقف حركة نقل الأبقار بسبب اشتباه في مرض التهاب الجلد العقدي الوطن | رصد توجهت وزارة الزراعة والثروة الحيوانية بتعليمات عاجلة إلى المنسقين في قطاعات الزراعة والثروة الحيوانية بالبلديات، بناءً على التقارير الواردة حول اشتباه في إصابات بمرض التهاب الجلد العقدي في الأبقار. ونظرًا لنتائج العينات التي أظهرت إصابة بالمرض الفيروسي، طلبت الوزارة بشكل عاجل إبلاغ الجهات المختصة بوقف حركة نقل وتنقل الأبقار بين البلديات ومنع عرضها في أسواق المواشي حتى إشعار آخر. تأتي هذه الإجراءات تنفيذًا لقانون الوقاية من الأمراض الحيوانية المعدية، وتهدف إلى حماية الثروة الحيوانية ومنع انتشار المرض.'
$ws.Range("C3").Value = 'https://www.libyaakhbar.com/libya-news/2263865.html'

# D3: datetime value with custom number format (164 created then superseded by 165, matching source)
$ws.Range("D3").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D3").Value = 45239.7459375
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
